$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44386
$ws.Cells.Item(2, 9).Value = 'Primera'
$ws.Cells.Item(2, 10).Value = 40
$ws.Cells.Item(2, 11).Value = 7000
$ws.Cells.Item(2, 12).Value = 7000
$ws.Cells.Item(2, 13).Value = 7000
$ws.Cells.Item(2, 15).Value = 'Región del Maule'
$ws.Cells.Item(2, 16).Value = 438

$ws.Cells.Item(3, 4).Value = 44396
$ws.Cells.Item(3, 9).Value = 'Primera'
$ws.Cells.Item(3, 10).Value = 80
$ws.Cells.Item(3, 11).Value = 7000
$ws.Cells.Item(3, 12).Value = 7000
$ws.Cells.Item(3, 13).Value = 7000
$ws.Cells.Item(3, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(3, 16).Value = 438

$ws.Cells.Item(4, 4).Value = 44355
$ws.Cells.Item(4, 9).Value = 'Primera'
$ws.Cells.Item(4, 10).Value = 30
$ws.Cells.Item(4, 11).Value = 8000
$ws.Cells.Item(4, 12).Value = 8000
$ws.Cells.Item(4, 13).Value = 8000
$ws.Cells.Item(4, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(4, 16).Value = 500

$ws.Cells.Item(5, 4).Value = 44348
$ws.Cells.Item(5, 9).Value = 'Primera'
$ws.Cells.Item(5, 10).Value = 35
$ws.Cells.Item(5, 11).Value = 7000
$ws.Cells.Item(5, 12).Value = 7000
$ws.Cells.Item(5, 13).Value = 7000
$ws.Cells.Item(5, 15).Value = 'Región del Maule'
$ws.Cells.Item(5, 16).Value = 438

$ws.Cells.Item(6, 4).Value = 44398
$ws.Cells.Item(6, 9).Value = 'Primera'
$ws.Cells.Item(6, 10).Value = 80
$ws.Cells.Item(6, 11).Value = 7000
$ws.Cells.Item(6, 12).Value = 7000
$ws.Cells.Item(6, 13).Value = 7000
$ws.Cells.Item(6, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(6, 16).Value = 438

$ws.Cells.Item(7, 4).Value = 44314
$ws.Cells.Item(7, 9).Value = 'Segunda'
$ws.Cells.Item(7, 10).Value = 20
$ws.Cells.Item(7, 11).Value = 5000
$ws.Cells.Item(7, 12).Value = 5000
$ws.Cells.Item(7, 13).Value = 5000
$ws.Cells.Item(7, 15).Value = 'Región del Maule'
$ws.Cells.Item(7, 16).Value = 312

$ws.Cells.Item(8, 4).Value = 44385
$ws.Cells.Item(8, 9).Value = 'Primera'
$ws.Cells.Item(8, 10).Value = 100
$ws.Cells.Item(8, 11).Value = 7000
$ws.Cells.Item(8, 12).Value = 7000
$ws.Cells.Item(8, 13).Value = 7000
$ws.Cells.Item(8, 15).Value = 'Región del Maule'
$ws.Cells.Item(8, 16).Value = 438

$ws.Cells.Item(9, 4).Value = 44362
$ws.Cells.Item(9, 9).Value = 'Primera'
$ws.Cells.Item(9, 10).Value = 25
$ws.Cells.Item(9, 11).Value = 8000
$ws.Cells.Item(9, 12).Value = 8000
$ws.Cells.Item(9, 13).Value = 8000
$ws.Cells.Item(9, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(9, 16).Value = 500

$ws.Cells.Item(10, 4).Value = 44315
$ws.Cells.Item(10, 9).Value = 'Primera'
$ws.Cells.Item(10, 10).Value = 40
$ws.Cells.Item(10, 11).Value = 7000
$ws.Cells.Item(10, 12).Value = 7000
$ws.Cells.Item(10, 13).Value = 7000
$ws.Cells.Item(10, 15).Value = 'Región del Maule'
$ws.Cells.Item(10, 16).Value = 438

$ws.Cells.Item(11, 4).Value = 44403
$ws.Cells.Item(11, 9).Value = 'Primera'
$ws.Cells.Item(11, 10).Value = 35
$ws.Cells.Item(11, 11).Value = 5000
$ws.Cells.Item(11, 12).Value = 5000
$ws.Cells.Item(11, 13).Value = 5000
$ws.Cells.Item(11, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(11, 16).Value = 312

$ws.Cells.Item(12, 4).Value = 44369
$ws.Cells.Item(12, 9).Value = 'Primera'
$ws.Cells.Item(12, 10).Value = 60
$ws.Cells.Item(12, 11).Value = 7000
$ws.Cells.Item(12, 12).Value = 7000
$ws.Cells.Item(12, 13).Value = 7000
$ws.Cells.Item(12, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(12, 16).Value = 438

$ws.Cells.Item(13, 4).Value = 44420
$ws.Cells.Item(13, 9).Value = 'Primera'
$ws.Cells.Item(13, 10).Value = 45
$ws.Cells.Item(13, 11).Value = 8000
$ws.Cells.Item(13, 12).Value = 8000
$ws.Cells.Item(13, 13).Value = 8000
$ws.Cells.Item(13, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(13, 16).Value = 500

$ws.Cells.Item(14, 4).Value = 44399
$ws.Cells.Item(14, 9).Value = 'Primera'
$ws.Cells.Item(14, 10).Value = 80
$ws.Cells.Item(14, 11).Value = 7000
$ws.Cells.Item(14, 12).Value = 7000
$ws.Cells.Item(14, 13).Value = 7000
$ws.Cells.Item(14, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(14, 16).Value = 438

$ws.Cells.Item(15, 4).Value = 44371
$ws.Cells.Item(15, 9).Value = 'Primera'
$ws.Cells.Item(15, 10).Value = 200
$ws.Cells.Item(15, 11).Value = 7000
$ws.Cells.Item(15, 12).Value = 7000
$ws.Cells.Item(15, 13).Value = 7000
$ws.Cells.Item(15, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(15, 16).Value = 438

$ws.Cells.Item(16, 4).Value = 44397
$ws.Cells.Item(16, 9).Value = 'Primera'
$ws.Cells.Item(16, 10).Value = 40
$ws.Cells.Item(16, 11).Value = 8000
$ws.Cells.Item(16, 12).Value = 8000
$ws.Cells.Item(16, 13).Value = 8000
$ws.Cells.Item(16, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(16, 16).Value = 500

$ws.Cells.Item(17, 4).Value = 44372
$ws.Cells.Item(17, 9).Value = 'Primera'
$ws.Cells.Item(17, 10).Value = 50
$ws.Cells.Item(17, 11).Value = 6000
$ws.Cells.Item(17, 12).Value = 7000
$ws.Cells.Item(17, 13).Value = 6400
$ws.Cells.Item(17, 15).Value = 'Región del Maule'
$ws.Cells.Item(17, 16).Value = 400

$ws.Cells.Item(18, 4).Value = 44308
$ws.Cells.Item(18, 9).Value = 'Primera'
$ws.Cells.Item(18, 10).Value = 75
$ws.Cells.Item(18, 11).Value = 5000
$ws.Cells.Item(18, 12).Value = 5000
$ws.Cells.Item(18, 13).Value = 5000
$ws.Cells.Item(18, 15).Value = 'Región del Maule'
$ws.Cells.Item(18, 16).Value = 312

$ws.Cells.Item(19, 4).Value = 44467
$ws.Cells.Item(19, 9).Value = 'Primera'
$ws.Cells.Item(19, 10).Value = 40
$ws.Cells.Item(19, 11).Value = 7000
$ws.Cells.Item(19, 12).Value = 7000
$ws.Cells.Item(19, 13).Value = 7000
$ws.Cells.Item(19, 15).Value = 'Región del Maule'
$ws.Cells.Item(19, 16).Value = 438

$ws.Cells.Item(20, 4).Value = 44313
$ws.Cells.Item(20, 9).Value = 'Primera'
$ws.Cells.Item(20, 10).Value = 20
$ws.Cells.Item(20, 11).Value = 7000
$ws.Cells.Item(20, 12).Value = 7000
$ws.Cells.Item(20, 13).Value = 7000
$ws.Cells.Item(20, 15).Value = 'Región del Maule'
$ws.Cells.Item(20, 16).Value = 438

$ws.Cells.Item(21, 4).Value = 44389
$ws.Cells.Item(21, 9).Value = 'Primera'
$ws.Cells.Item(21, 10).Value = 55
$ws.Cells.Item(21, 11).Value = 7000
$ws.Cells.Item(21, 12).Value = 7000
$ws.Cells.Item(21, 13).Value = 7000
$ws.Cells.Item(21, 15).Value = 'Región del Maule'
$ws.Cells.Item(21, 16).Value = 438

$ws.Cells.Item(22, 4).Value = 44305
$ws.Cells.Item(22, 9).Value = 'Primera'
$ws.Cells.Item(22, 10).Value = 35
$ws.Cells.Item(22, 11).Value = 7000
$ws.Cells.Item(22, 12).Value = 7000
$ws.Cells.Item(22, 13).Value = 7000
$ws.Cells.Item(22, 15).Value = 'Región del Maule'
$ws.Cells.Item(22, 16).Value = 438

$ws.Cells.Item(23, 4).Value = 44392
$ws.Cells.Item(23, 9).Value = 'Primera'
$ws.Cells.Item(23, 10).Value = 95
$ws.Cells.Item(23, 11).Value = 7000
$ws.Cells.Item(23, 12).Value = 7000
$ws.Cells.Item(23, 13).Value = 7000
$ws.Cells.Item(23, 15).Value = 'Región del Maule'
$ws.Cells.Item(23, 16).Value = 438

$ws.Cells.Item(24, 4).Value = 44354
$ws.Cells.Item(24, 9).Value = 'Primera'
$ws.Cells.Item(24, 10).Value = 100
$ws.Cells.Item(24, 11).Value = 8000
$ws.Cells.Item(24, 12).Value = 9000
$ws.Cells.Item(24, 13).Value = 8500
$ws.Cells.Item(24, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(24, 16).Value = 531

$ws.Cells.Item(25, 4).Value = 44354
$ws.Cells.Item(25, 9).Value = 'Primera'
$ws.Cells.Item(25, 10).Value = 80
$ws.Cells.Item(25, 11).Value = 9000
$ws.Cells.Item(25, 12).Value = 9000
$ws.Cells.Item(25, 13).Value = 9000
$ws.Cells.Item(25, 15).Value = 'Región del Maule'
$ws.Cells.Item(25, 16).Value = 562

$ws.Cells.Item(26, 4).Value = 44312
$ws.Cells.Item(26, 9).Value = 'Primera'
$ws.Cells.Item(26, 10).Value = 40
$ws.Cells.Item(26, 11).Value = 7000
$ws.Cells.Item(26, 12).Value = 7000
$ws.Cells.Item(26, 13).Value = 7000
$ws.Cells.Item(26, 15).Value = 'Región del Maule'
$ws.Cells.Item(26, 16).Value = 438
